$d = $word.ActiveDocument

# --- Paragraph "Button search": move the _GoBack bookmark to sit right
# after the "search" run (before the closing spellEnd proofErr) ---
$pSearch = $d.Paragraphs.Item(39)
$rSearch = $pSearch.Range
if ($rSearch.Text.TrimEnd() -ne "Button search") {
    throw "Unexpected paragraph 39 content: [$($rSearch.Text)]"
}
$xmlSearch = @'
<w:p w14:paraId="6FCA25CE" w14:textId="05987594" w:rsidR="00303863" w:rsidRDefault="00303863" w:rsidP="00303863" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Button </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>search</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p>
'@
$rSearch.InsertXML($xmlSearch)

# --- Paragraph "Grid patientinformation" (under PatientView, the last
# list item) -> wrap the phrase in parentheses: "(Grid patientinformation)".
# Word's proofer now also wraps "patientinformation" with spellStart/spellEnd.
# The old trailing _GoBack bookmark is dropped here since it moved above. ---
$pGrid = $d.Paragraphs.Item(41)
$rGrid = $pGrid.Range
if ($rGrid.Text.TrimEnd() -ne "Grid patientinformation") {
    throw "Unexpected paragraph 41 content: [$($rGrid.Text)]"
}
$xmlGrid = @'
<w:p w14:paraId="0341645E" w14:textId="03DBC0A6" w:rsidR="00303863" w:rsidRDefault="00303863" w:rsidP="00303863" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Grid</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>patientinformation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@
$rGrid.InsertXML($xmlGrid)
